# Harmonize parameters w Leander thesis
$wb = $excel.ActiveWorkbook

# LOHC sheet: Spec capex trailer (euros) 150000 -> 660000
$wsLOHC = $wb.Worksheets.Item("LOHC")
$wsLOHC.Range("B12").Value = 660000

# NH3 sheet: Costs for driver (euros/h) 20 -> 2.85
$wsNH3 = $wb.Worksheets.Item("NH3")
$wsNH3.Range("B5").Value = 2.85

# NH3 sheet: Truck lifetime (a) 12 -> 8
$wsNH3.Range("B11").Value = 8

# NH3 sheet: Spec capex trailer (euros) 190000 -> 210000
$wsNH3.Range("B12").Value = 210000
